# Update "想去人数" (want-to-go count) figures for several events,
# matching the scraped data refresh described in the commit message.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F9").Value  = 2108
$ws1.Range("F13").Value = 1009
$ws1.Range("F14").Value = 154
$ws1.Range("F15").Value = 2157
$ws1.Range("F16").Value = 612
$ws1.Range("F17").Value = 11197
$ws1.Range("F18").Value = 1133

# --- Sheet "全部类型" (aggregated view of all categories) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F15").Value = 2108
$ws4.Range("F21").Value = 1009
$ws4.Range("F23").Value = 154
$ws4.Range("F26").Value = 2157
$ws4.Range("F27").Value = 612
$ws4.Range("F28").Value = 11197
$ws4.Range("F31").Value = 1133
